$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-9 from 2023-10-25 (45224)
# to 2023-11-03 (45233), keeping existing cell formatting/style intact.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45233
}
